$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "307.44"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-4.11%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "39.85"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-6.11%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.119"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-1.44%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07732"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.235"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-1.81%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.614"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-10.97%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8934"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-4.31%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1002"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-9.44%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1739"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-6.53%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09005"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-4.27%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04441"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-5.82%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1054"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.40%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001256"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-4.05%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005845"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.94%"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "2,410.60%"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.11%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.419"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-4.49%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3318"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-2.10%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.054"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-4.95%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1348"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-3.07%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2759"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "8.22%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04139"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.24%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001207"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-3.02%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004063"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-5.45%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001301"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "8.27%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02349"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-13.05%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05202"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-6.22%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007921"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-2.29%"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-5.45%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.006281"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-4.25%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.001951"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-6.57%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008211"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-0.70%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3324"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-4.60%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006506"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-6.16%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.06%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "98.14%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.003505"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "4.34%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.06%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002001"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.06%"
